$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the tenant_id / tenant_id_lbl column (column H) entirely from the
# "base/org" excel template: both the header/comment row (row 1) and the
# data template row (row 2) reference tenant_id_lbl only in this column,
# so deleting the whole column removes both shared-string entries and
# shifts the subsequent columns (I, J) left to (H, I).
$ws.Columns("H").Delete()
